# B6-PowerPoint.pptx edit — commit "Wed, Jul 15, 2020  7:07:07 PM"
#
# 1) Three tables (on slides 14, 15 and 16) switch from the custom
#    "Table_0" style to the built-in PowerPoint table style
#    {24DAAA9C-1FA9-45F2-9CDA-E6BB38859268}.
# 2) The deck's theme colour palette is switched from the "Integral /
#    Red Violet" palette to the standard Office palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------
$tableSlideIndexes = 14, 15, 16
$newTableStyleId = "{24DAAA9C-1FA9-45F2-9CDA-E6BB38859268}"

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the theme colour palette for the Office defaults -------
$officeColors = @{
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $themeColors.Item($idx).RGB = $officeColors[$idx]
}
